# The deck ships two DrawingML theme parts:
#   theme1.xml -> "Office Theme" (default Office colors)
#   theme2.xml -> "Integral"     (the theme actually applied to the
#                                  slide master / whole presentation)
# The authored edit swaps the two themes' contents, so the presentation's
# applied design becomes the plain "Office Theme" colors, while the
# "Integral" palette moves to the other (unused) theme slot.
#
# font scheme and format scheme (fills/lines/effects) are identical between
# the two themes - only the colour scheme differs - so re-pointing the
# live presentation theme's 12 colour slots to the Office Theme palette
# reproduces the applied visual change.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.Theme.ThemeColorScheme

# Office Theme colour scheme (was theme1.xml), expressed as VBA RGB()
# (low byte = R, mid byte = G, high byte = B) integers so they land in
# ppt/theme/theme2.xml, the theme actually referenced by the presentation.
$cs.Item(1).RGB  = 0         # dk1      000000
$cs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388   # dk2      44546A
$cs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501   # accent2  ED7D31
$cs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Item(8).RGB  = 49407     # accent4  FFC000
$cs.Item(9).RGB  = 12874308  # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
